$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID Venta"
$ws.Range("B1").Value = "Fecha"
$ws.Range("C1").Value = "Hora"
$ws.Range("D1").Value = "Vendedor"
$ws.Range("E1").Value = "Productos"
$ws.Range("F1").Value = "Total"
$ws.Range("G1").Value = "Personas"

$ws.Range("A2").Value = "V-1769891068199"
$ws.Range("B2").Value = "31/1/2026"
$ws.Range("C2").Value = "03:24 p. m."
$ws.Range("D2").Value = "Martha"
$ws.Range("E2").Value = "Aguardiente Amarillo Botella (x1), Aguardiente Amarillo Media (x1), Ron 5 años Botella (x1), Cerveza Corona (x1)"
$ws.Range("F2").Value = 274000
$ws.Range("G2").Value = 0
